$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = 1.29
$ws.Range("Q2").Value = 1.63
$ws.Range("AG2").Value = 1000

# Row 3 updates
$ws.Range("T3").Value = 1.84
$ws.Range("U3").Value = 1.94
$ws.Range("Z3").Value = 11.5
